# Auto-generated edit script: updates cryptos list values per commit
# "Updated cryptos list on Mon Apr 10 22:23:10 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple text/string cell updates (Coin name / Link / Price-with-dual-dots / Volume%) ---
$simpleUpdates = @{
    'D2' = '29.587.81'
    'E2' = '  +3.17%  '
    'D3' = '1.905.66'
    'E3' = '  +1.18%  '
    'E4' = '  -0.78%  '
    'E5' = '  -0.58%  '
    'E6' = '  -0.77%  '
    'E7' = '  +0.63%  '
    'E8' = '  +0.42%  '
    'E9' = '  +0.49%  '
    'B10' = 'OKB'
    'C10' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'E10' = '  +1.85%  '
    'B11' = 'Polygon'
    'C11' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'E11' = '  +0.14%  '
    'E12' = '  +0.00%  '
    'D13' = '1.906.54'
    'E13' = '  +1.02%  '
    'E14' = '  +1.13%  '
    'E15' = '  +0.33%  '
    'E16' = '  -0.94%  '
    'E17' = '  +1.73%  '
    'E18' = '  +0.11%  '
    'E19' = '  +0.03%  '
    'E20' = '  +0.88%  '
    'E21' = '  -0.75%  '
    'E22' = '  +0.74%  '
    'D23' = '29.577.23'
    'E23' = '  +3.04%  '
    'E24' = '  +0.09%  '
    'E25' = '  -1.84%  '
    'D26' = '2.116.15'
    'E26' = '  +0.68%  '
    'B27' = 'Monero'
    'C27' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'E27' = '  -1.34%  '
    'B28' = 'EthereumClassic'
    'C28' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'E28' = '  +0.72%  '
    'E29' = '  +2.07%  '
    'E30' = '  +0.81%  '
    'E31' = '  +1.35%  '
    'E32' = '  -0.27%  '
    'E33' = '  +6.05%  '
    'E34' = '  +0.53%  '
    'E35' = '  +0.81%  '
    'E36' = '  +0.96%  '
    'B37' = 'Algorand'
    'C37' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'E37' = '  +0.11%  '
    'B38' = 'FraxShare'
    'C38' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'E38' = '  +1.05%  '
    'E39' = '  +2.41%  '
    'E40' = '  +2.56%  '
    'E41' = '  +0.88%  '
    'E42' = '  -2.53%  '
    'E43' = '  +0.80%  '
    'E44' = '  +0.05%  '
    'E45' = '  +1.00%  '
    'E46' = '  -0.86%  '
    'E47' = '  +1.13%  '
    'E48' = '  +0.81%  '
    'E49' = '  +0.98%  '
    'E50' = '  -3.09%  '
    'E51' = '  +0.70%  '
}
foreach ($addr in $simpleUpdates.Keys) {
    $ws.Range($addr).Value = $simpleUpdates[$addr]
}

# --- Price cells whose new text would otherwise be auto-coerced to a number by Excel ---
# (e.g. "1.002", "42.50") -- force as text via NumberFormat "@", then restore the
# original cell style so no residual formatting/style change is left behind.
$textPreservedUpdates = @{
    'D4' = '1.002'
    'D5' = '314.71'
    'D6' = '1.001'
    'D7' = '0.5155'
    'D8' = '0.3976'
    'D9' = '0.08484'
    'D10' = '42.50'
    'D11' = '1.119'
    'D12' = '6.312'
    'D14' = '20.79'
    'D15' = '7.327'
    'D16' = '1.002'
    'D17' = '93.21'
    'D18' = '0.00001110'
    'D19' = '0.06748'
    'D20' = '17.96'
    'D21' = '1.001'
    'D22' = '6.034'
    'D25' = '2.213'
    'D27' = '159.21'
    'D28' = '20.93'
    'D30' = '128.26'
    'D31' = '1.067'
    'D32' = '0.1053'
    'D33' = '6.176'
    'D34' = '3.649'
    'D35' = '0.02493'
    'D36' = '0.06610'
    'D37' = '0.2203'
    'D38' = '9.065'
    'D39' = '5.238'
    'D40' = '1.234'
    'D41' = '0.6536'
    'D43' = '11.32'
    'D44' = '0.6093'
    'D45' = '13.19'
    'D49' = '123.97'
    'D51' = '78.00'
}
foreach ($addr in $textPreservedUpdates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $textPreservedUpdates[$addr]
    $cell.Style = $origStyle
}

